$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.878.54"
$ws.Range("E2").Value = "  -4.35%  "
$ws.Range("D3").Value = "3.801.88"
$ws.Range("E3").Value = "  -4.99%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'588.05"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("D6").Value = "'164.00"
$ws.Range("E6").Value = "  +2.28%  "
$ws.Range("D7").Value = "'0.659"
$ws.Range("E7").Value = "  -3.78%  "
$ws.Range("E8").Value = "  +0.29%  "
$ws.Range("D9").Value = "'0.734"
$ws.Range("E9").Value = "  -2.39%  "
$ws.Range("E10").Value = "  +2.19%  "
$ws.Range("D11").Value = "'52.04"
$ws.Range("E11").Value = "  -3.89%  "
$ws.Range("E12").Value = "  -1.07%  "
$ws.Range("D13").Value = "'11.05"
$ws.Range("E13").Value = "  +0.41%  "
$ws.Range("D14").Value = "4.411.29"
$ws.Range("E14").Value = "  -4.86%  "
$ws.Range("D15").Value = "3.820.79"
$ws.Range("E15").Value = "  -4.51%  "
$ws.Range("D16").Value = "'20.57"
$ws.Range("E16").Value = "  +0.64%  "
$ws.Range("D17").Value = "'13.59"
$ws.Range("E17").Value = "  -4.15%  "
$ws.Range("E18").Value = "  -6.88%  "
$ws.Range("D20").Value = "69.854.46"
$ws.Range("E20").Value = "  -4.06%  "
$ws.Range("D21").Value = "'433.38"
$ws.Range("E21").Value = "  -0.40%  "
$ws.Range("D22").Value = "'4.65"
$ws.Range("E22").Value = "  -3.30%  "
$ws.Range("D23").Value = "'92.46"
$ws.Range("E23").Value = "  -4.23%  "
$ws.Range("E24").Value = "  -6.60%  "
$ws.Range("D25").Value = "'13.68"
$ws.Range("E25").Value = "  -4.11%  "
$ws.Range("D26").Value = "'11.06"
$ws.Range("E26").Value = "  -2.26%  "
$ws.Range("E27").Value = "  -11.71%  "
$ws.Range("D28").Value = "'5.94"
$ws.Range("E28").Value = "  -0.47%  "
$ws.Range("D29").Value = "'10.30"
$ws.Range("E29").Value = "  -1.88%  "
$ws.Range("D31").Value = "'7.94"
$ws.Range("E31").Value = "  +0.19%  "
$ws.Range("E32").Value = "  -4.59%  "
$ws.Range("D33").Value = "'47.58"
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("E34").Value = "  -6.38%  "
$ws.Range("D35").Value = "'68.36"
$ws.Range("E35").Value = "  -3.28%  "
$ws.Range("D36").Value = "'0.0₃0964"
$ws.Range("E36").Value = "  +9.23%  "
$ws.Range("D37").Value = "'622.38"
$ws.Range("E37").Value = "  -7.47%  "
$ws.Range("D38").Value = "'0.418"
$ws.Range("E38").Value = "  -5.56%  "
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("E40").Value = "  -0.06%  "
$ws.Range("E41").Value = "  -2.81%  "
$ws.Range("D42").Value = "'3.20"
$ws.Range("E42").Value = "  -4.54%  "
$ws.Range("E43").Value = "  +17.85%  "
$ws.Range("E44").Value = "  -5.82%  "
$ws.Range("E45").Value = "  +4.82%  "
$ws.Range("D46").Value = "'9.76"
$ws.Range("E46").Value = "  -9.16%  "
$ws.Range("E47").Value = "  -5.80%  "
$ws.Range("E48").Value = "  -15.76%  "
$ws.Range("D49").Value = "'3.25"
$ws.Range("E49").Value = "  -5.09%  "
$ws.Range("D50").Value = "2.819.33"
$ws.Range("E50").Value = "  -1.21%  "
$ws.Range("E51").Value = "  -0.26%  "
